$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the full D2:E51 range to Text format so numeric-looking strings
# (e.g. "1.002", "28.177.77") are preserved as text, not converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

# Apply the updated values from the diff (Coin name, Link, Price, Volume columns).
$ws.Range("D2").Value = '28.177.77'
$ws.Range("E2").Value = '  +0.72%  '
$ws.Range("D3").Value = '1.803.09'
$ws.Range("E3").Value = '  +3.44%  '
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  -0.37%  '
$ws.Range("D5").Value = '336.36'
$ws.Range("E5").Value = '  +0.25%  '
$ws.Range("E6").Value = '  -0.12%  '
$ws.Range("D7").Value = '0.4649'
$ws.Range("E7").Value = '  +23.43%  '
$ws.Range("D8").Value = '0.3700'
$ws.Range("E8").Value = '  +10.61%  '
$ws.Range("D9").Value = '45.48'
$ws.Range("E9").Value = '  +0.09%  '
$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").Value = '0.07670'
$ws.Range("E10").Value = '  +6.67%  '
$ws.Range("B11").Value = 'Polygon'
$ws.Range("C11").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D11").Value = '1.155'
$ws.Range("E11").Value = '  +4.00%  '
$ws.Range("D12").Value = '22.60'
$ws.Range("E12").Value = '  +1.05%  '
$ws.Range("D13").Value = '1.001'
$ws.Range("E13").Value = '  -0.32%  '
$ws.Range("D14").Value = '6.374'
$ws.Range("E14").Value = '  +3.88%  '
$ws.Range("D15").Value = '7.402'
$ws.Range("E15").Value = '  +3.90%  '
$ws.Range("D16").Value = '1.799.64'
$ws.Range("E16").Value = '  +2.80%  '
$ws.Range("D17").Value = '0.00001098'
$ws.Range("E17").Value = '  +4.03%  '
$ws.Range("D18").Value = '0.06738'
$ws.Range("E18").Value = '  +2.49%  '
$ws.Range("D19").Value = '82.89'
$ws.Range("E19").Value = '  +4.05%  '
$ws.Range("D20").Value = '1.000'
$ws.Range("E20").Value = '  -0.13%  '
$ws.Range("D21").Value = '17.46'
$ws.Range("E21").Value = '  +3.74%  '
$ws.Range("D22").Value = '6.437'
$ws.Range("E22").Value = '  +3.34%  '
$ws.Range("D23").Value = '28.146.30'
$ws.Range("E23").Value = '  +0.58%  '
$ws.Range("D24").Value = '11.93'
$ws.Range("E24").Value = '  +2.42%  '
$ws.Range("D25").Value = '2.416'
$ws.Range("E25").Value = '  +0.79%  '
$ws.Range("D26").Value = '20.81'
$ws.Range("E26").Value = '  +5.25%  '
$ws.Range("D27").Value = '2.399'
$ws.Range("E27").Value = '  +3.52%  '
$ws.Range("D28").Value = '152.49'
$ws.Range("E28").Value = '  -0.95%  '
$ws.Range("D29").Value = '2.004.13'
$ws.Range("E29").Value = '  +2.75%  '
$ws.Range("D30").Value = '134.94'
$ws.Range("E30").Value = '  +2.60%  '
$ws.Range("D31").Value = '1.269'
$ws.Range("E31").Value = '  +1.84%  '
$ws.Range("D32").Value = '4.047'
$ws.Range("E32").Value = '  +0.58%  '
$ws.Range("D33").Value = '0.09639'
$ws.Range("E33").Value = '  +10.55%  '
$ws.Range("D34").Value = '5.913'
$ws.Range("E34").Value = '  +2.50%  '
$ws.Range("D35").Value = '0.2252'
$ws.Range("E35").Value = '  +7.06%  '
$ws.Range("B36").Value = 'Aptos'
$ws.Range("C36").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D36").Value = '12.25'
$ws.Range("E36").Value = '  +0.51%  '
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").Value = '0.02380'
$ws.Range("E37").Value = '  +2.78%  '
$ws.Range("D38").Value = '0.06403'
$ws.Range("E38").Value = '  +3.50%  '
$ws.Range("D39").Value = '0.6733'
$ws.Range("E39").Value = '  +0.92%  '
$ws.Range("D40").Value = '5.284'
$ws.Range("E40").Value = '  +2.67%  '
$ws.Range("D41").Value = '1.525'
$ws.Range("E41").Value = '  +5.58%  '
$ws.Range("D42").Value = '1.236'
$ws.Range("E42").Value = '  +2.01%  '
$ws.Range("D43").Value = '8.144'
$ws.Range("E43").Value = '  +2.67%  '
$ws.Range("D44").Value = '14.13'
$ws.Range("E44").Value = '  +3.03%  '
$ws.Range("E45").Value = '  -0.17%  '
$ws.Range("D46").Value = '0.6182'
$ws.Range("E46").Value = '  +2.49%  '
$ws.Range("D47").Value = '3.843'
$ws.Range("E47").Value = '  +0.53%  '
$ws.Range("D48").Value = '130.34'
$ws.Range("E48").Value = '  +1.84%  '
$ws.Range("D49").Value = '2.070'
$ws.Range("E49").Value = '  +2.77%  '
$ws.Range("D50").Value = '1.187'
$ws.Range("E50").Value = '  +1.26%  '
$ws.Range("D51").Value = '0.07140'
$ws.Range("E51").Value = '  +0.23%  '

# Clear the temporary text formatting we applied so the cell style index
# matches the original (unstyled) cells.
$ws.Range("D2:E51").Style = "Normal"

